# Add the new monthly data row (01-10-2021) to Sheet1, right after the
# existing last row (46), mirroring the layout of the prior rows:
# col A = period label (text), cols B:F = numeric index values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 47

# Column A holds the period label as literal text (e.g. "01-09-2021" in the
# row above), not an Excel date. Force text formatting before assigning the
# value so Excel's automatic date recognition doesn't convert the string
# into a date serial number, then clear the formatting override again so
# the cell ends up styled the same as its neighbours (no explicit style).
$cellA = $ws.Cells.Item($row, 1)
$cellA.NumberFormat = "@"
$cellA.Value = "01-10-2021"
$cellA.ClearFormats()

$ws.Cells.Item($row, 2).Value = 112.94
$ws.Cells.Item($row, 3).Value = 110.94
$ws.Cells.Item($row, 4).Value = 114.84
$ws.Cells.Item($row, 5).Value = 110.7
$ws.Cells.Item($row, 6).Value = 123.91
